# Update "Attendance Roster_Splunk Certified Architect_October 06-10_VC00529848.xlsx"
# Adds per-learner VM access details (IP1 / IP2 / Username=ubuntu) next to the
# existing Username/Password credential table on Sheet2.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(2)

# ---------------------------------------------------------------------------
# 0. Pre-format the new cells (F:I, rows 11:23) by copying the plain row
#    formatting already used by the C column in each row (rows 12-23 use the
#    plain/unfilled look, row 11 is the bold/filled header look carried by
#    D11). Formatting first, values after - matches how the cells ended up
#    styled in the saved workbook.
# ---------------------------------------------------------------------------
$ws.Range("C12").Copy() | Out-Null
$ws.Range("F11:F23").PasteSpecial(-4122) | Out-Null   # xlPasteFormats
$ws.Range("G12:I23").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

$ws.Range("D11").Copy() | Out-Null
$ws.Range("G11:I11").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = 0

# ---------------------------------------------------------------------------
# 1. Values - entered in the same order they were originally authored:
#    IP1 column (G13:G23) pasted first, then IP2 column (H12:H23), then the
#    G12 value typed in, then the two header labels, then the Username
#    ("ubuntu") column, then the Username header label.
# ---------------------------------------------------------------------------
$ip1 = @{
    13 = "174.129.73.140"
    14 = "52.91.55.48"
    15 = "34.229.155.235"
    16 = "98.80.9.40"
    17 = "34.202.230.107"
    18 = "52.55.213.186"
    19 = "54.242.108.128"
    20 = "34.236.146.139"
    21 = "34.229.57.121"
    22 = "98.84.127.232"
    23 = "54.221.6.74"
}
foreach ($r in 13..23) {
    $ws.Cells.Item($r, 7).Value = $ip1[$r]
}

$ip2 = @{
    12 = "34.228.155.29"
    13 = "3.208.12.191"
    14 = "13.220.50.93"
    15 = "100.26.111.168"
    16 = "18.208.220.175"
    17 = "34.224.93.22"
    18 = "54.227.96.203"
    19 = "23.20.122.36"
    20 = "3.91.66.211"
    21 = "3.94.115.94"
    22 = "34.234.63.205"
    23 = "54.90.72.181"
}
foreach ($r in 12..23) {
    $ws.Cells.Item($r, 8).Value = $ip2[$r]
}

$ws.Range("G12").Value = "54.167.22.31  "

$ws.Range("G11").Value = "IP1"
$ws.Range("H11").Value = "IP2"

foreach ($r in 12..23) {
    $ws.Cells.Item($r, 9).Value = "ubuntu"
}

$ws.Range("I11").Value = "Username"

# ---------------------------------------------------------------------------
# 2. New column widths for G, H, I (bestFit-style, sized to the new content).
# ---------------------------------------------------------------------------
$ws.Columns.Item(7).ColumnWidth = 13.0
$ws.Columns.Item(8).ColumnWidth = 37.5
$ws.Columns.Item(9).ColumnWidth = 13.666666666666666

# ---------------------------------------------------------------------------
# 3. Update the active selection to match the saved view state.
# ---------------------------------------------------------------------------
$ws.Activate()
$ws.Range("J9:K9").Select() | Out-Null
